# Delete rows 2 through 7 (years 2000年, 2005年-2009年), which shifts
# the remaining rows (2010年-2013年) up to become rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:C7").EntireRow.Delete()
